$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) keeps its values as text, matching the sheet's existing
# string-typed data (many prices contain multiple "." separators and would
# otherwise be auto-coerced into numbers by Excel's input parser).
$ws.Range("D2:D51").NumberFormat = "@"

# --- Coin / Link swaps (rows 44/45 and 50/51 traded ranking positions) ---
$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"

# --- Price (D) and Volume(1h) (E) refresh ---
$ws.Range("D2").Value = "43.032.15"
$ws.Range("E2").Value = "  +0.60%  "
$ws.Range("D3").Value = "2.302.57"
$ws.Range("E3").Value = "  -0.42%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").Value = "310.14"
$ws.Range("E5").Value = "  -2.64%  "
$ws.Range("D6").Value = "104.43"
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -0.46%  "
$ws.Range("D10").Value = "39.54"
$ws.Range("E10").Value = "  -1.54%  "
$ws.Range("E11").Value = "  -0.65%  "
$ws.Range("E12").Value = "  -3.32%  "
$ws.Range("E13").Value = "  -0.10%  "
$ws.Range("D14").Value = "0.992"
$ws.Range("E14").Value = "  +1.47%  "
$ws.Range("D15").Value = "2.781.45"
$ws.Range("E15").Value = "  +4.42%  "
$ws.Range("D16").Value = "15.34"
$ws.Range("E16").Value = "  -0.51%  "
$ws.Range("D17").Value = "2.303.20"
$ws.Range("E17").Value = "  -1.34%  "
$ws.Range("D18").Value = "42.823.84"
$ws.Range("E18").Value = "  +0.22%  "
$ws.Range("D19").Value = "7.32"
$ws.Range("E19").Value = "  -2.68%  "
$ws.Range("E20").Value = "  -1.32%  "
$ws.Range("D21").Value = "13.41"
$ws.Range("E21").Value = "  +0.83%  "
$ws.Range("D22").Value = "73.36"
$ws.Range("E22").Value = "  -0.63%  "
$ws.Range("D23").Value = "3.44"
$ws.Range("E23").Value = "  -3.12%  "
$ws.Range("D24").Value = "268.22"
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("D25").Value = "2.21"
$ws.Range("E25").Value = "  -2.13%  "
$ws.Range("E26").Value = "  +0.65%  "
$ws.Range("D27").Value = "10.90"
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("D28").Value = "7.31"
$ws.Range("E28").Value = "  +17.53%  "
$ws.Range("D29").Value = "2.29"
$ws.Range("E29").Value = "  -1.35%  "
$ws.Range("D30").Value = "22.28"
$ws.Range("E30").Value = "  -1.86%  "
$ws.Range("D31").Value = "36.25"
$ws.Range("E31").Value = "  -4.46%  "
$ws.Range("D32").Value = "164.92"
$ws.Range("E32").Value = "  -0.52%  "
$ws.Range("E33").Value = "  -3.73%  "
$ws.Range("E34").Value = "  -1.23%  "
$ws.Range("E35").Value = "  +2.15%  "
$ws.Range("E36").Value = "  -3.72%  "
$ws.Range("D37").Value = "4.54"
$ws.Range("E37").Value = "  -1.84%  "
$ws.Range("E38").Value = "  -1.45%  "
$ws.Range("E39").Value = "  +2.07%  "
$ws.Range("E40").Value = "  -2.49%  "
$ws.Range("D41").Value = "109.18"
$ws.Range("E41").Value = "  +10.42%  "
$ws.Range("E42").Value = "  -1.37%  "
$ws.Range("D43").Value = "70.80"
$ws.Range("E43").Value = "  +0.83%  "
$ws.Range("D44").Value = "0.226"
$ws.Range("E44").Value = "  +0.20%  "
$ws.Range("D45").Value = "1.01"
$ws.Range("E45").Value = "  +0.24%  "
$ws.Range("D46").Value = "12.31"
$ws.Range("E46").Value = "  -0.71%  "
$ws.Range("D47").Value = "1.730.06"
$ws.Range("E47").Value = "  +6.75%  "
$ws.Range("D48").Value = "110.86"
$ws.Range("E48").Value = "  -3.62%  "
$ws.Range("D49").Value = "77.80"
$ws.Range("E49").Value = "  -5.39%  "
$ws.Range("D50").Value = "8.64"
$ws.Range("E50").Value = "  -3.02%  "
$ws.Range("D51").Value = "5.13"
$ws.Range("E51").Value = "  -3.11%  "
